# Deploy the implementation guide: refresh the generated IG metadata
# (status + timestamp) and make sure the "top + wrap" alignment that the
# header/data cell styles already declare is actually turned on.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: Status -> draft, Date -> new generation timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B6").Value = "draft"
$meta.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# --- Turn on the vertical-top / wrap-text alignment for every sheet ---
# The header row and the data rows already carry an <alignment> of
# vertical="top" wrapText="true" in the stylesheet, but it was never
# flagged as "applied", so Excel was ignoring it. Re-asserting the same
# alignment on every used cell flips that flag on for the whole workbook.
foreach ($sheet in $wb.Worksheets) {
    $used = $sheet.UsedRange
    $used.VerticalAlignment = -4160   # xlVAlignTop
    $used.WrapText = $true
}
